$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "29.843.82"
Set-TextValue $ws "E2" "  -1.12%  "
Set-TextValue $ws "D3" "1.892.59"
Set-TextValue $ws "E3" "  -0.98%  "
Set-TextValue $ws "E4" "  +0.06%  "
Set-TextValue $ws "D5" "0.7830"
Set-TextValue $ws "E5" "  -4.85%  "
Set-TextValue $ws "D6" "243.95"
Set-TextValue $ws "E6" "  +0.03%  "
Set-TextValue $ws "E7" "  +0.03%  "
Set-TextValue $ws "D8" "0.3143"
Set-TextValue $ws "E8" "  -3.56%  "
Set-TextValue $ws "D9" "25.29"
Set-TextValue $ws "E9" "  -5.91%  "
Set-TextValue $ws "D10" "0.07194"
Set-TextValue $ws "E10" "  +1.97%  "
Set-TextValue $ws "D11" "0.08091"
Set-TextValue $ws "E11" "  -0.14%  "
Set-TextValue $ws "D12" "0.7640"
Set-TextValue $ws "E12" "  -1.28%  "
Set-TextValue $ws "D13" "5.496"
Set-TextValue $ws "E13" "  +3.75%  "
Set-TextValue $ws "D14" "1.861.23"
Set-TextValue $ws "E14" "  -2.66%  "
Set-TextValue $ws "D15" "92.31"
Set-TextValue $ws "E15" "  -1.23%  "
Set-TextValue $ws "D16" "6.151"
Set-TextValue $ws "E16" "  +3.90%  "
Set-TextValue $ws "D17" "29.848.09"
Set-TextValue $ws "E17" "  -1.09%  "
Set-TextValue $ws "D18" "13.97"
Set-TextValue $ws "E18" "  -1.89%  "
Set-TextValue $ws "D19" "243.35"
Set-TextValue $ws "E19" "  -1.25%  "
Set-TextValue $ws "D20" "0.000007786"
Set-TextValue $ws "E20" "  -0.23%  "
Set-TextValue $ws "D21" "1.002"
Set-TextValue $ws "E21" "  -0.02%  "
Set-TextValue $ws "D22" "2.146.17"
Set-TextValue $ws "E22" "  -0.96%  "
Set-TextValue $ws "D23" "8.128"
Set-TextValue $ws "E23" "  +14.47%  "
Set-TextValue $ws "E24" "  +0.11%  "
Set-TextValue $ws "D25" "0.1639"
Set-TextValue $ws "E25" "  -2.30%  "
Set-TextValue $ws "D26" "9.423"
Set-TextValue $ws "E26" "  +0.90%  "
Set-TextValue $ws "D27" "163.03"
Set-TextValue $ws "E27" "  -2.47%  "
Set-TextValue $ws "E28" "  -1.40%  "
Set-TextValue $ws "D29" "2.051"
Set-TextValue $ws "E29" "  -2.99%  "
Set-TextValue $ws "D30" "1.412"
Set-TextValue $ws "E30" "  +2.92%  "
Set-TextValue $ws "D31" "1.549"
Set-TextValue $ws "E31" "  +1.35%  "
Set-TextValue $ws "D32" "4.495"
Set-TextValue $ws "E32" "  +4.25%  "
Set-TextValue $ws "D33" "4.110"
Set-TextValue $ws "E33" "  +0.07%  "
Set-TextValue $ws "D34" "0.05560"
Set-TextValue $ws "D35" "1.268"
Set-TextValue $ws "E35" "  -0.67%  "
Set-TextValue $ws "D36" "0.7439"
Set-TextValue $ws "D37" "0.9981"
Set-TextValue $ws "E37" "  -0.20%  "
Set-TextValue $ws "D38" "2.618"
Set-TextValue $ws "E38" "  -2.84%  "
Set-TextValue $ws "D39" "0.01921"
Set-TextValue $ws "E39" "  -0.31%  "
Set-TextValue $ws "D40" "2.786"
Set-TextValue $ws "E40" "  -0.28%  "
Set-TextValue $ws "D41" "1.147.57"
Set-TextValue $ws "E41" "  +14.06%  "
Set-TextValue $ws "D42" "73.71"
Set-TextValue $ws "E42" "  +0.29%  "
Set-TextValue $ws "D43" "0.4417"
Set-TextValue $ws "E43" "  -1.17%  "
Set-TextValue $ws "E44" "  -1.91%  "
Set-TextValue $ws "D45" "0.8509"
Set-TextValue $ws "E45" "  -0.28%  "
Set-TextValue $ws "E46" "  +0.01%  "
Set-TextValue $ws "D47" "103.92"
Set-TextValue $ws "E47" "  +1.14%  "
Set-TextValue $ws "D48" "1.877"
Set-TextValue $ws "E48" "  -1.97%  "
Set-TextValue $ws "D49" "9.979"
Set-TextValue $ws "E49" "  +1.25%  "
Set-TextValue $ws "D50" "7.456"
Set-TextValue $ws "E50" "  -1.89%  "
Set-TextValue $ws "D51" "2.997"
Set-TextValue $ws "E51" "  +9.73%  "
